$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("N:N").Insert()

$ws.Range("N4").Value = "Mã gói dịch vụ"
